# "further cleaning to metadata"
#
# - s2cDNAProtocol value changes from E7760 -> E7420 for every data row
#   (G2:G49), and picks up a dedicated Arial/11/black font (previously it
#   shared the plain Arial/11 font used elsewhere on the sheet).
# - roboticS2Prep (H2:H49) becomes a real "=FALSE()" formula instead of a
#   bare boolean literal (still displays/evaluates to FALSE).
# - The sheet's live selection follows the edited column: G2:G49 (was H2:H49).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 49

$protocolRange = $ws.Range("G$firstRow`:G$lastRow")

# New protocol code.
$protocolRange.Value = "E7420"

# Give the protocol column its own distinct font (Arial 11, black) rather
# than reusing the generic Arial/11 style.
$protocolRange.Font.Name = "Arial"
$protocolRange.Font.Size = 11
$protocolRange.Font.Color = 0

# roboticS2Prep becomes a live formula cell (=FALSE()) on every row instead
# of a stored boolean literal. Looping row-by-row (rather than assigning the
# whole range at once) keeps each cell an independent formula.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Formula = "=FALSE()"
}

# Move the active selection to the column that was just edited.
$protocolRange.Select() | Out-Null
